$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.265.64"
$ws.Range("D3").Value = "3.938.39"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "493.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("E7").Value = "  -1.06%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.177"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000352"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.43%  "
$ws.Range("D14").Value = "4.565.68"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").Value = "3.930.87"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.64%  "
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("E18").Value = "  +4.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("D20").Value = "69.257.83"
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "439.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.74%  "
$ws.Range("E22").Value = "  -0.40%  "
$ws.Range("E23").Value = "  -2.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.17%  "
$ws.Range("E26").Value = "  +5.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "703.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.56%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.77%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.130"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.467"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +16.68%  "
$ws.Range("E35").Value = "  -2.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "61.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.30%  "
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "40.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.66%  "
$ws.Range("E39").Value = "  +0.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.997"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.82%  "
$ws.Range("E44").Value = "  -4.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.41%  "
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.19%  "
$ws.Range("E48").Value = "  +6.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.33%  "
$ws.Range("D50").Value = "0.0₆0347"
$ws.Range("E50").Value = "  -5.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "144.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.89%  "
